$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create row 10 with the same formatting as the existing table rows (centered
# alignment, as used by rows 2-9) before filling in its values.
$ws.Range("A9:I9").Copy()
$ws.Range("A10:I10").PasteSpecial(-4122)

# Row 8: Router G0/0/0.10 (replaces the former placeholder "Switch" row)
$ws.Range("A8").Value = "Router G0/0/0.10"
$ws.Range("B8").Value = "192.168.10.1"
$ws.Range("C8").Value = "/24"
$ws.Range("D8").Value = "192.168.10.0"
$ws.Range("E8").Value = "LAN Applicatif"
$ws.Range("F8").Value = "-"
$ws.Range("G8").Value = "-"
$ws.Range("H8").Value = "-"
$ws.Range("I8").Value = "G0/0/0"

# Row 9: Router G0/0/0.20 (replaces the former placeholder "Router" row)
$ws.Range("A9").Value = "Router G0/0/0.20"
$ws.Range("B9").Value = "192.168.20.1"
$ws.Range("C9").Value = "/24"
$ws.Range("D9").Value = "192.168.20.0"
$ws.Range("E9").Value = "LAN Client"
$ws.Range("F9").Value = "-"
$ws.Range("G9").Value = "-"
$ws.Range("H9").Value = "-"
$ws.Range("I9").Value = "G0/0/0"

# Row 10: Router G0/0/0.30 (new row)
$ws.Range("A10").Value = "Router G0/0/0.30"
$ws.Range("B10").Value = "192.168.30.1"
$ws.Range("C10").Value = "/24"
$ws.Range("D10").Value = "192.168.30.0"
$ws.Range("E10").Value = "LAN DMZ"
$ws.Range("F10").Value = "-"
$ws.Range("G10").Value = "-"
$ws.Range("H10").Value = "-"
$ws.Range("I10").Value = "G0/0/0"

# Match the final selection state shown in the saved workbook
$ws.Range("I10").Select()
